$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Rename "Sheet2" -> "Financial_Institution"
# ---------------------------------------------------------------------------
$fi = $wb.Worksheets.Item("Sheet2")
$fi.Name = "Financial_Institution"

# ---------------------------------------------------------------------------
# Populate header row 1 (columns A-L, then N-Q; M1 intentionally left blank)
# ---------------------------------------------------------------------------
$fi.Range("A1").Value = "Institution_Name"
$fi.Range("B1").Value = "FI_Number"
$fi.Range("C1").Value = "Address_Line1"
$fi.Range("D1").Value = "Address_Line2"
$fi.Range("E1").Value = "City"
$fi.Range("F1").Value = "State"
$fi.Range("G1").Value = "ZIP"
$fi.Range("H1").Value = "ZIP_SUFFIX"
$fi.Range("I1").Value = "Corporate_Structure"
$fi.Range("J1").Value = "TAX_Identification_Number"
$fi.Range("K1").Value = "OmniBus_Account_Number"
$fi.Range("L1").Value = "Special_Instructions"
$fi.Range("N1").Value = "Omni_Serve"
$fi.Range("O1").Value = "NSCC_ID"
$fi.Range("P1").Value = "SuccessMessage_FailureMessage"
$fi.Range("Q1").Value = "MessageKeyword"

# ---------------------------------------------------------------------------
# Populate message-keyword lookup values in column Q (rows 2-5)
# ---------------------------------------------------------------------------
$fi.Range("Q2").Value = "AccountExist"
$fi.Range("Q3").Value = "WrongTaxIdNumber"
$fi.Range("Q4").Value = "WrongNSCCID"
$fi.Range("Q5").Value = "Success"

# ---------------------------------------------------------------------------
# Touch rows 6-9 (A:Q) so they materialize as used (empty) rows, matching
# the extended sheet dimension A1:Q9 that Excel wrote after formatting.
# ---------------------------------------------------------------------------
$fi.Range("A6:Q9").ClearFormats()

# ---------------------------------------------------------------------------
# Column widths (best-fit in the source workbook). The underlying engine
# quantizes ColumnWidth to 1/6-character steps, so the values below are the
# closest settable approximations of the recorded bestFit widths.
# ---------------------------------------------------------------------------
$fi.Columns.Item(1).ColumnWidth = 15.833333333333334
$fi.Columns.Item(2).ColumnWidth = 10.0
$fi.Columns.Item(3).ColumnWidth = 13.166666666666666
$fi.Columns.Item(4).ColumnWidth = 13.166666666666666
$fi.Columns.Item(8).ColumnWidth = 9.833333333333334
$fi.Columns.Item(9).ColumnWidth = 18.5
$fi.Columns.Item(10).ColumnWidth = 25.5
$fi.Columns.Item(11).ColumnWidth = 25.166666666666668
$fi.Columns.Item(14).ColumnWidth = 11.0
$fi.Columns.Item(16).ColumnWidth = 30.166666666666668
$fi.Columns.Item(17).ColumnWidth = 15.833333333333334

# ---------------------------------------------------------------------------
# View state: Financial_Institution becomes the active/selected tab, with
# Q5 as the active cell (matching the recorded sheetView selection).
# ---------------------------------------------------------------------------
$fi.Activate()
$fi.Range("Q5").Select()
